# Update the cached "datetimeFigureOut" Date placeholder text from
# 9/15/2023 to 11/27/2023 across every slide layout, the slide master,
# and the notes master (the footer/date placeholders that PowerPoint
# keeps in sync when the deck's header & footer date is refreshed).

$p = $ppt.ActivePresentation
$newDate = "11/27/2023"

function Update-DatePlaceholders($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*") {
            $shp.TextFrame.TextRange.Text = $newDate
        }
    }
}

# Slide master
Update-DatePlaceholders $p.SlideMaster.Shapes

# Every slide layout (custom layout) hanging off the slide master
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-DatePlaceholders $layouts.Item($li).Shapes
}

# Notes master (its Date placeholder only persists edits made through
# the HeadersFooters object, unlike the slide master/layouts above)
$p.NotesMaster.HeadersFooters.DateAndTime.Text = $newDate
